$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 / 40: Hedera and VeChain swap places (name + link) ---
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

# --- Price column (D): force text entry so values like "1.006" are not
#     re-interpreted as numbers, then clear the temporary formatting so the
#     cell style index is left untouched (matches original "General" style). ---
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"
$ws.Range("D2").Value = "24.330.71"
$ws.Range("D3").Value = "1.647.33"
$ws.Range("D4").Value = "1.006"
$ws.Range("D5").Value = "310.62"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.3645"
$ws.Range("D8").Value = "46.71"
$ws.Range("D9").Value = "0.3227"
$ws.Range("D10").Value = "1.113"
$ws.Range("D11").Value = "0.06981"
$ws.Range("D13").Value = "5.905"
$ws.Range("D14").Value = "19.22"
$ws.Range("D15").Value = "6.538"
$ws.Range("D16").Value = "1.644.51"
$ws.Range("D17").Value = "0.00001029"
$ws.Range("D18").Value = "0.06579"
$ws.Range("D20").Value = "77.59"
$ws.Range("D21").Value = "5.884"
$ws.Range("D22").Value = "15.46"
$ws.Range("D23").Value = "12.39"
$ws.Range("D24").Value = "24.339.80"
$ws.Range("D25").Value = "2.483"
$ws.Range("D26").Value = "2.275"
$ws.Range("D27").Value = "145.38"
$ws.Range("D28").Value = "18.41"
$ws.Range("D29").Value = "1.827.86"
$ws.Range("D30").Value = "123.22"
$ws.Range("D31").Value = "1.160"
$ws.Range("D32").Value = "4.052"
$ws.Range("D33").Value = "5.579"
$ws.Range("D34").Value = "0.08406"
$ws.Range("D35").Value = "1.661"
$ws.Range("D36").Value = "11.98"
$ws.Range("D37").Value = "5.131"
$ws.Range("D38").Value = "1.234"
$ws.Range("D39").Value = "0.02201"
$ws.Range("D40").Value = "0.05938"
$ws.Range("D41").Value = "0.2036"
$ws.Range("D42").Value = "8.041"
$ws.Range("D44").Value = "0.5825"
$ws.Range("D45").Value = "3.750"
$ws.Range("D46").Value = "12.41"
$ws.Range("D47").Value = "0.5546"
$ws.Range("D48").Value = "121.39"
$ws.Range("D49").Value = "1.930"
$ws.Range("D50").Value = "0.06871"
$ws.Range("D51").Value = "1.172"
$dRange.ClearFormats()

# --- Volume(1h) column (E): plain text percentages, safe to set directly. ---
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("E8").Value = "  -5.85%  "
$ws.Range("E9").Value = "  -6.51%  "
$ws.Range("E10").Value = "  -8.22%  "
$ws.Range("E11").Value = "  -7.57%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  -6.77%  "
$ws.Range("E14").Value = "  -9.52%  "
$ws.Range("E15").Value = "  -7.67%  "
$ws.Range("E16").Value = "  -3.85%  "
$ws.Range("E17").Value = "  -9.34%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -8.91%  "
$ws.Range("E21").Value = "  -8.26%  "
$ws.Range("E22").Value = "  -11.05%  "
$ws.Range("E23").Value = "  -7.03%  "
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -18.99%  "
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -10.34%  "
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("E30").Value = "  -7.42%  "
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("E33").Value = "  -19.68%  "
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("E35").Value = "  -7.48%  "
$ws.Range("E36").Value = "  -14.11%  "
$ws.Range("E37").Value = "  -8.89%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -8.97%  "
$ws.Range("E40").Value = "  -11.23%  "
$ws.Range("E41").Value = "  -9.16%  "
$ws.Range("E42").Value = "  -14.40%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -9.87%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("E46").Value = "  -11.28%  "
$ws.Range("E47").Value = "  -10.30%  "
$ws.Range("E48").Value = "  -6.96%  "
$ws.Range("E49").Value = "  -9.75%  "
$ws.Range("E50").Value = "  -6.20%  "
